# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 30 (pushing the existing rows 30-197
# down to 31-198) and populate it with the new week's price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 30. Excel will shift rows 30..197
# down to 31..198 and copy the formatting of the row above into the new row.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new record's values.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R have the same value on every data row of
# this sheet, so fill them in explicitly along with the week's new figures.
$ws.Range("A30").Value = 8
$ws.Range("B30").Value = "Terminal La Palmera de La Serena"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44550
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 100112012
$ws.Range("G30").Value = "Espinaca"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 2000
$ws.Range("K30").Value = 400
$ws.Range("L30").Value = 500
$ws.Range("M30").Value = 450
$ws.Range("N30").Value = "$/atado 300 a 500 gramos"
$ws.Range("O30").Value = "Provincia del Elquí"
$ws.Range("P30").Value = 900
$ws.Range("Q30").Value = 0.5
$ws.Range("R30").Value = "Hortaliza"
